$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.10989254549611616
$ws.Range("B1").Value = 0.10927061526663806
$ws.Range("A2").Value = -0.10977651752371642
$ws.Range("B2").Value = 0.10828946798977679
$ws.Range("A3").Value = -0.058573684777384116
$ws.Range("B3").Value = 0.05814255628773779
$ws.Range("A4").Value = -0.050142556402933636
$ws.Range("B4").Value = 0.049756872769643579
$ws.Range("A5").Value = -0.04675687282628882
$ws.Range("B5").Value = 0.045449513510352979
$ws.Range("A6").Value = -0.012989811738657764
$ws.Range("B6").Value = 0.012796032979709437
$ws.Range("A7").Value = -0.0027960331283747308
$ws.Range("B7").Value = 0.0027614351856404262
$ws.Range("A8").Value = 0.0072385646660020342
$ws.Range("B8").Value = -0.0072774980671521661
$ws.Range("A9").Value = 0.0092774980213587988
$ws.Range("B9").Value = -0.0093036451320616109
$ws.Range("A10").Value = 0.011303645090297465
$ws.Range("B10").Value = -0.011303112097333923
$ws.Range("A11").Value = 0.014303112043469568
$ws.Range("B11").Value = -0.014304099590471964
$ws.Range("A12").Value = 0.01780409953173967
$ws.Range("B12").Value = -0.017830444956672409
$ws.Range("A13").Value = 0.021330444905847123
$ws.Range("B13").Value = -0.021351547836917284
$ws.Range("A14").Value = 0.029351547733848626
$ws.Range("B14").Value = -0.029387291137120464
$ws.Range("A15").Value = 0.030387291124641003
$ws.Range("B15").Value = -0.030434118036518854
$ws.Range("A16").Value = 0.032434118016994695
$ws.Range("B16").Value = -0.032630861951517431
$ws.Range("A17").Value = -0.0040029393653471246
$ws.Range("B17").Value = 0.0039999999668793862
$ws.Range("A18").Value = 0.0065059021816402662
$ws.Range("B18").Value = -0.0067290917203521872
$ws.Range("A19").Value = 0.0039176768863029565
$ws.Range("B19").Value = -0.0045074558928965835
$ws.Range("A20").Value = 0.0050966581035787328
$ws.Range("B20").Value = -0.0051776744325433555
$ws.Range("A21").Value = -0.0040056686080687243
$ws.Range("B21").Value = 0.0039999999494808591
$ws.Range("A22").Value = -0.045715783274657085
$ws.Range("B22").Value = 0.045501954189040461
$ws.Range("A23").Value = -0.04050195426784331
$ws.Range("B23").Value = 0.040099400145074959
$ws.Range("A24").Value = -0.020099400416586199
$ws.Range("B24").Value = 0.019999999724596762
$ws.Range("A25").Value = -0.030459701823037477
$ws.Range("B25").Value = 0.030405410166020275
$ws.Range("A26").Value = -0.027905410224125404
$ws.Range("B26").Value = 0.027837843297405129
$ws.Range("A27").Value = -0.025337843356938894
$ws.Range("B27").Value = 0.024949469040839389
$ws.Range("A28").Value = -0.022949469099677877
$ws.Range("B28").Value = 0.022703848216295341
$ws.Range("A29").Value = -0.015703848340516302
$ws.Range("B29").Value = 0.015645633172836071
$ws.Range("A30").Value = 0.044354366046847726
$ws.Range("B30").Value = -0.04456960870572102
$ws.Range("A31").Value = 0.051569608591027105
$ws.Range("B31").Value = -0.051643670924033103
$ws.Range("A32").Value = -0.0040008195577776462
$ws.Range("B32").Value = 0.003999999925984099
